$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H: header "Save" in row 1 (same style as the other header cells),
# and the save-flag values for rows 2-6.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 1
